$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Total" header in G1, matching the bold/centered style used by the
# other header cells (B1:F1)
$ws.Range("G1").Value = "Total"
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").HorizontalAlignment = -4108

# Fill in the total counts for each realm row
$ws.Range("G2").Value = 7
$ws.Range("G3").Value = 38
$ws.Range("G4").Value = 16
$ws.Range("G5").Value = 28
$ws.Range("G6").Value = 18
